# Updates the order records on the "CheckOutRegression" sheet used by the
# checkout tests (order id / subtotal / tax / total for the two sample
# orders), refreshing them with the current test order numbers.
#
# Note: the saved window position (xWindow) is a cosmetic, last-used-window
# artifact that Excel's object model does not expose for writing, so it is
# left untouched here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CheckOutRegression")

# --- Row 2: order "00001255" -> "00001269" (amount updated) ---
$ws.Range("L2").Value = "ORDER NUMBER`n00001269"
$ws.Range("N2").Value = "£77.16"
$ws.Range("P2").Value = "£10.99"
$ws.Range("O2").Value = "Your order includes £12.86 tax."
$ws.Range("M2").Value = "£88.15"

# --- Row 5: order "00001257" -> "00001271" (amounts unchanged) ---
$ws.Range("L5").Value = "ORDER NUMBER`n00001271"
$ws.Range("N5").Value = "£53.36"
$ws.Range("P5").Value = "£10.99"
$ws.Range("O5").Value = "Your order includes £8.89 tax."
$ws.Range("M5").Value = "£64.35"

# The orderTax column (O) no longer needs to be as wide once the longer,
# best-fit text is replaced, so give it an explicit (narrower) width.
$ws.Columns.Item(15).ColumnWidth = 16.666666666
